# Rename item class names to follow the new naming convention.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Weapon_MeleeAttackBasic"      # was BasicAttackItem (Strike)
$ws.Range("B4").Value = "Weapon_RangedAttackFireBall"  # was BasicRangedAttackItem (Fireball)
$ws.Range("B5").Value = "Weapon_RangedAttackBow"       # was Weapon_BasicRangedBow (Swift Shot)

$wb.Save()
